# Daily attendance processing - 2026-01-18 23:58:50
# Reorders the comma-separated "Recorded By" email lists in column G
# for several rows on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value  = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"
$ws.Range("G3").Value  = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System"
$ws.Range("G4").Value  = "asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G5").Value  = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G6").Value  = "alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G7").Value  = "lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G13").Value = "esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G17").Value = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
$ws.Range("G19").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G20").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
